$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumValue($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

function Set-BoolValue($row, $col, $bool) {
    $ws.Cells.Item($row, $col).Value = $bool
}

# The "Startdatum"/"Slutdatum" columns (Y, AA) hold plain text that looks
# like a date (e.g. "2023-09-19"). Typing that directly into a General cell
# would make Excel auto-convert it into a date serial number, which is not
# what the source data contains. Temporarily mark those cells as Text before
# entering the values, then restore the cells to the default "Normal" style
# once the literal text is safely stored (so no lingering text-number-format
# is left behind on the cells).
$tempTextCells = @()
foreach ($r in 25..27) {
    foreach ($c in 25, 27) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $tempTextCells += $cell
    }
}

# Columns I (9), K (11), AT (46) and AY (51) hold an explicit empty-string
# value in the source data (as opposed to being completely absent from the
# row). Force each such cell to exist by briefly formatting it as Text and
# assigning an empty string; the temporary Text formatting is reverted,
# along with the date cells above, at the end of the script.
foreach ($r in 25..27) {
    foreach ($c in 9, 11, 46, 51) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = ""
        $tempTextCells += $cell
    }
}

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$r = 25
Set-NumValue  $r 1  112196967               # A  Id
Set-NumValue  $r 2  43467                   # B  Taxonsorteringsordning
$ws.Cells.Item($r,3).Value  = "Ovaliderad"  # C  Valideringsstatus
$ws.Cells.Item($r,4).Value  = "LC"          # D  Rodlistade
Set-NumValue  $r 5  101735                  # E  TaxonId
$ws.Cells.Item($r,6).Value  = "Jättesvampmal"      # F  Artnamn
$ws.Cells.Item($r,7).Value  = "Scardia boletella"  # G  Vetenskapligt namn
$ws.Cells.Item($r,8).Value  = "(Fabricius, 1794)"  # H  Auktor
$ws.Cells.Item($r,16).Value = "Flugmötesskogen , Srm"  # P  Lokalnamn
Set-NumValue  $r 17 580549.8394260339       # Q  Ost
Set-NumValue  $r 18 6579319.7617336         # R  Nord
Set-NumValue  $r 19 5                       # S  Noggrannhet
$ws.Cells.Item($r,20).Value = "Södermanland" # T  Lan
$ws.Cells.Item($r,21).Value = "Eskilstuna"   # U  Kommun
$ws.Cells.Item($r,22).Value = "Södermanland" # V  Provins
$ws.Cells.Item($r,23).Value = "Eskilstuna"   # W  Forsamling
$ws.Cells.Item($r,25).Value = "2023-09-19"            # Y  Startdatum
$ws.Cells.Item($r,26).Value = "12:30"        # Z  Starttid
$ws.Cells.Item($r,27).Value = "2023-09-19"            # AA Slutdatum
$ws.Cells.Item($r,28).Value = "12:30"        # AB Sluttid
Set-BoolValue $r 30 $false                   # AD Ej aterfunnen
Set-BoolValue $r 31 $false                   # AE Osaker artbestamning
Set-BoolValue $r 33 $false                   # AG Ospontan
$ws.Cells.Item($r,49).Value = "Ella Axelsson Elfving"  # AW Rapportor
$ws.Cells.Item($r,50).Value = "Ella Axelsson Elfving"  # AX Observatorer

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$r = 26
Set-NumValue  $r 1  112195278
Set-NumValue  $r 2  8377
$ws.Cells.Item($r,3).Value  = "Ovaliderad"
$ws.Cells.Item($r,4).Value  = "LC"
Set-NumValue  $r 5  106545
$ws.Cells.Item($r,6).Value  = "Mindre märgborre"
$ws.Cells.Item($r,7).Value  = "Tomicus minor"
$ws.Cells.Item($r,8).Value  = "(Hartig, 1834)"
$ws.Cells.Item($r,16).Value = "Flugmötesskogen , Srm"
Set-NumValue  $r 17 580549.8394260339
Set-NumValue  $r 18 6579319.7617336
Set-NumValue  $r 19 5
$ws.Cells.Item($r,20).Value = "Södermanland"
$ws.Cells.Item($r,21).Value = "Eskilstuna"
$ws.Cells.Item($r,22).Value = "Södermanland"
$ws.Cells.Item($r,23).Value = "Eskilstuna"
$ws.Cells.Item($r,25).Value = "2023-09-19"
$ws.Cells.Item($r,26).Value = "11:00"
$ws.Cells.Item($r,27).Value = "2023-09-19"
$ws.Cells.Item($r,28).Value = "11:00"
Set-BoolValue $r 30 $false
Set-BoolValue $r 31 $false
Set-BoolValue $r 33 $false
$ws.Cells.Item($r,49).Value = "Ella Axelsson Elfving"
$ws.Cells.Item($r,50).Value = "Ella Axelsson Elfving"

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$r = 27
Set-NumValue  $r 1  112194720
Set-NumValue  $r 2  56414
$ws.Cells.Item($r,3).Value  = "Ovaliderad"
$ws.Cells.Item($r,4).Value  = "NT"
Set-NumValue  $r 5  100049
$ws.Cells.Item($r,6).Value  = "Spillkråka"
$ws.Cells.Item($r,7).Value  = "Dryocopus martius"
$ws.Cells.Item($r,8).Value  = "(Linnaeus, 1758)"
$ws.Cells.Item($r,16).Value = "Flugmötesskogen , Srm"
Set-NumValue  $r 17 580549.8394260339
Set-NumValue  $r 18 6579319.7617336
Set-NumValue  $r 19 5
$ws.Cells.Item($r,20).Value = "Södermanland"
$ws.Cells.Item($r,21).Value = "Eskilstuna"
$ws.Cells.Item($r,22).Value = "Södermanland"
$ws.Cells.Item($r,23).Value = "Eskilstuna"
$ws.Cells.Item($r,25).Value = "2023-09-19"
$ws.Cells.Item($r,26).Value = "10:56"
$ws.Cells.Item($r,27).Value = "2023-09-19"
$ws.Cells.Item($r,28).Value = "10:56"
$ws.Cells.Item($r,29).Value = "Hördes i omgivningen långa rop"  # AC Publik kommentar
Set-BoolValue $r 30 $false
Set-BoolValue $r 31 $false
Set-BoolValue $r 33 $false
$ws.Cells.Item($r,49).Value = "Ella Axelsson Elfving"
$ws.Cells.Item($r,50).Value = "Ella Axelsson Elfving"

# ---------------------------------------------------------------------------
# Restore every cell that was temporarily forced to Text format (the Y/AA
# date cells and the explicit-empty-string cells) back to the default
# "Normal" cell style now that their values are safely stored.
# ---------------------------------------------------------------------------
foreach ($cell in $tempTextCells) {
    $cell.Style = "Normal"
}
